# Insert a new weekly price record as row 140, pushing the existing
# rows 140..201 down to 141..202 (dimension grows from A1:R201 to A1:R202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 140.
$ws.Rows("140:140").Insert()

# The row that used to be 140 is now 141; clone its (shared, unchanged)
# columns into the new row 140 so formatting/labels match, then overwrite
# the handful of cells that hold the new record's own data.
$src = $ws.Range("A141:R141").Value2
$ws.Range("A140:R140").Value2 = $src

# New record's own values (Fecha, Volumen, Precio minimo/maximo/promedio,
# Precio $/Kg). Calidad stays "Primera" as copied from the row above.
$ws.Cells.Item(140, 4).Value2  = 45007   # D140 Fecha
$ws.Cells.Item(140, 10).Value2 = 500     # J140 Volumen
$ws.Cells.Item(140, 11).Value2 = 1000    # K140 Precio minimo
$ws.Cells.Item(140, 12).Value2 = 1000    # L140 Precio maximo
$ws.Cells.Item(140, 13).Value2 = 1000    # M140 Precio promedio ponderado
$ws.Cells.Item(140, 16).Value2 = 1000    # P140 Precio $/Kg
